$d = $word.ActiveDocument

# Update the Script Execution Order value for GONet.GONetGlobal from
# -200 to the correct value of -32000.
$d.Content.Find.Execute(
    "at a value of -200 and",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "at a value of -32000 and", 2
)
